$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-10-07 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-08 Tuesday", 2)

$tbl = $d.Tables.Item(1)

$values = @(
  @("78×32=", "79×67=", "88×86=", "45×75=", "55×87="),
  @("34×47=", "76×47=", "85×30=", "90×49=", "30×81="),
  @("29×18=", "44×64=", "53×32=", "81×45=", "73×44="),
  @("45×99=", "79×72=", "93×29=", "33×90=", "20×19="),
  @("77×13=", "29×63=", "11×65=", "93×19=", "25×11=")
)

$rows = @(1, 5, 10, 15, 20)

for ($r = 0; $r -lt 5; $r++) {
  $rowIndex = $rows[$r]
  for ($c = 1; $c -le 5; $c++) {
    $cell = $tbl.Cell($rowIndex, $c)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $values[$r][$c - 1]
  }
}
